$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.917.14"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "3.065.10"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "527.08"
$ws.Range("E5").Value = "  +6.45%  "
$ws.Range("D6").Value = "143.73"
$ws.Range("E6").Value = "  +6.76%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +5.89%  "
$ws.Range("D9").Value = "7.68"
$ws.Range("E9").Value = "  +6.83%  "
$ws.Range("E10").Value = "  +8.24%  "
$ws.Range("E11").Value = "  +6.33%  "
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "3.593.34"
$ws.Range("E13").Value = "  +3.16%  "
$ws.Range("D14").Value = "27.49"
$ws.Range("E14").Value = "  +9.53%  "
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  +17.58%  "
$ws.Range("D16").Value = "57.948.95"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").Value = "6.23"
$ws.Range("E17").Value = "  +7.34%  "
$ws.Range("D18").Value = "3.060.00"
$ws.Range("E18").Value = "  +2.87%  "
$ws.Range("D19").Value = "13.25"
$ws.Range("E19").Value = "  +7.85%  "
$ws.Range("D20").Value = "8.20"
$ws.Range("E20").Value = "  +5.76%  "
$ws.Range("D21").Value = "341.55"
$ws.Range("E21").Value = "  +5.13%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "5.69"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("E24").Value = "  +7.54%  "
$ws.Range("D25").Value = "64.93"
$ws.Range("E25").Value = "  +5.61%  "
$ws.Range("E26").Value = "  +6.52%  "
$ws.Range("D27").Value = "0.0₃0981"
$ws.Range("E27").Value = "  +10.15%  "
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "6.98"
$ws.Range("E29").Value = "  +7.71%  "
$ws.Range("D30").Value = "7.43"
$ws.Range("E30").Value = "  +10.47%  "
$ws.Range("E31").Value = "  +7.34%  "
$ws.Range("E32").Value = "  +7.53%  "
$ws.Range("D33").Value = "21.11"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("D34").Value = "4.80"
$ws.Range("E34").Value = "  +8.34%  "
$ws.Range("D35").Value = "157.16"
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  +7.90%  "
$ws.Range("E37").Value = "  +4.82%  "
$ws.Range("D38").Value = "26.43"
$ws.Range("E38").Value = "  +14.66%  "
$ws.Range("D39").Value = "0.0706"
$ws.Range("E39").Value = "  +5.86%  "
$ws.Range("D40").Value = "3.102.35"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("E41").Value = "  +3.90%  "
$ws.Range("D42").Value = "3.92"
$ws.Range("E42").Value = "  +10.36%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.48"
$ws.Range("E43").Value = "  +6.32%  "
$ws.Range("D44").Value = "0.667"
$ws.Range("E44").Value = "  +4.67%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "2.338.83"
$ws.Range("E46").Value = "  +5.48%  "
$ws.Range("E47").Value = "  +3.93%  "
$ws.Range("D48").Value = "2.02"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("E49").Value = "  +5.80%  "
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("D51").Value = "20.24"
$ws.Range("E51").Value = "  +7.38%  "
